$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# --- Update values for row 2 (aggregate) ---
# Row 2
$ws.Range("D2").Value = -0.0818
$ws.Range("E2").Value = -0.102
$ws.Range("F2").Value = 0.02155
$ws.Range("G2").Value = 0.10865137667039
$ws.Range("H2").Value = 0.10865137667039
$ws.Range("I2").Value = 0.08488828736278817
$ws.Range("J2").Value = 0.07209348667685544
$ws.Range("K2").Value = 2957.57
$ws.Range("L2").Value = 0.04018950789910668
$ws.Range("M2").Value = 3072.412
$ws.Range("N2").Value = 0.03811726932343476
$ws.Range("O2").Value = 1.038829850181061
$ws.Range("P2").Value = 3032.73
$ws.Range("Q2").Value = 0.03762496247093824
$ws.Range("R2").Value = 1.025412754389583
$ws.Range("S2").Value = 39.68199999999991
$ws.Range("T2").Value = 0.01291558553996011
$ws.Range("U2").Value = 45693.6
$ws.Range("V2").Value = 0.5668885740445287
$ws.Range("W2").Value = 0.1185683796126239
$ws.Range("X2").Value = 0.07516690505239665
$ws.Range("Y2").Value = 0.04340147456022724
$ws.Range("Z2").Value = 1.731965762288543
$ws.Range("AA2").Value = 0.05728703015591918
$ws.Range("AB2").Value = 0.06091989149169285
$ws.Range("AC2").Value = -0.003580702335058996
$ws.Range("AD2").Value = 30634.81
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 30634.81
$ws.Range("AG2").Value = -15058.79
$ws.Range("AH2").Value = 0.2753962840913453
$ws.Range("AI2").Value = 0.4212266840612526
$ws.Range("AJ2").Value = -0.2297459120325893
$ws.Range("AK2").Value = -0.5570313098204817
$ws.Range("AL2").Value = 1119.99
$ws.Range("AM2").Value = 1119.99
$ws.Range("AN2").Value = 4.636144749705652
$ws.Range("AO2").Value = 5.577710515272457
$ws.Range("AP2").Value = -2.278934656210369
$ws.Range("AQ2").Value = 5.577710515272457

# Row 3
$ws.Range("B3").Value = "Just Group plc (LSE:JUST)"
$ws.Range("F3").Value = -0.0138
$ws.Range("G3").Value = 0.1734107935938271
$ws.Range("H3").Value = 0.1734107935938271
$ws.Range("I3").Value = 0.2079090215782154
$ws.Range("J3").Value = 0.1693494536664662
$ws.Range("K3").Value = 553
$ws.Range("L3").Value = 0.1240410928177291
$ws.Range("M3").Value = 33.7
$ws.Range("N3").Value = 0.03401635207429091
$ws.Range("O3").Value = 0.06094032549728753
$ws.Range("P3").Value = 33.7
$ws.Range("Q3").Value = 0.03401635207429091
$ws.Range("R3").Value = 0.06094032549728753
$ws.Range("U3").Value = 2086.4
$ws.Range("V3").Value = 2.105985666700313
$ws.Range("W3").Value = 0.2038935181771256
$ws.Range("X3").Value = 0.0958609313413455
$ws.Range("Y3").Value = 0.1080325868357801
$ws.Range("Z3").Value = 2.339525608732158
$ws.Range("AA3").Value = 0.3961973836774977
$ws.Range("AB3").Value = 0.06333563903058731
$ws.Range("AC3").Value = 0.3328617446469104
$ws.Range("AD3").Value = 845.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 845.5
$ws.Range("AG3").Value = -1240.9
$ws.Range("AH3").Value = 0.4604618233307918
$ws.Range("AI3").Value = 0.2108531384822564
$ws.Range("AJ3").Value = 4.959632294164668
$ws.Range("AK3").Value = -0.6451260722641019
$ws.Range("AL3").Value = 216.7
$ws.Range("AM3").Value = 216.7
$ws.Range("AN3").Value = 0.9103143841515935
$ws.Range("AO3").Value = 4.277341947392709
$ws.Range("AP3").Value = -1.336024978466839
$ws.Range("AQ3").Value = 4.277341947392709

# Row 4
$ws.Range("B4").Value = "Phoenix Group Holdings plc (LSE:PHNX)"
$ws.Range("D4").Value = 0.287
$ws.Range("E4").Value = 0.206
$ws.Range("F4").Value = -0.06619999999999999
$ws.Range("G4").Value = 0.06855738402707694
$ws.Range("H4").Value = 0.06855738402707694
$ws.Range("I4").Value = 0.09758408337280911
$ws.Range("J4").Value = 0.07374876899304557
$ws.Range("K4").Value = 685.9
$ws.Range("L4").Value = 0.04708943491305034
$ws.Range("M4").Value = 454.4
$ws.Range("N4").Value = 0.0475936108929039
$ws.Range("O4").Value = 0.6624872430383437
$ws.Range("P4").Value = 454.4
$ws.Range("Q4").Value = 0.0475936108929039
$ws.Range("R4").Value = 0.6624872430383437
$ws.Range("U4").Value = 8054.3
$ws.Range("V4").Value = 0.8436030374443572
$ws.Range("W4").Value = 0.09811887561690866
$ws.Range("X4").Value = 0.08051248950915564
$ws.Range("Y4").Value = 0.01760638610775302
$ws.Range("Z4").Value = 4.638526208521751
$ws.Range("AA4").Value = 0.3420855978204581
$ws.Range("AB4").Value = 0.0607589449641641
$ws.Range("AC4").Value = 0.281326652856294
$ws.Range("AD4").Value = 4613.4
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 4613.4
$ws.Range("AG4").Value = -3440.900000000001
$ws.Range("AH4").Value = 0.3257843781115607
$ws.Range("AI4").Value = 0.386673371888358
$ws.Range("AJ4").Value = -0.5634723086496579
$ws.Range("AK4").Value = -0.887584801506436
$ws.Range("AL4").Value = 234
$ws.Range("AM4").Value = 234
$ws.Range("AN4").Value = 3.170939583476527
$ws.Range("AO4").Value = 6.074358974358975
$ws.Range("AP4").Value = -2.365042270946457
$ws.Range("AQ4").Value = 6.074358974358975

# Row 5
$ws.Range("B5").Value = "Chesnara plc (LSE:CSN)"
$ws.Range("D5").Value = -0.0818
$ws.Range("E5").Value = -0.102
$ws.Range("G5").Value = 0.2860295740454645
$ws.Range("H5").Value = 0.2860295740454645
$ws.Range("I5").Value = 0.09335687486206135
$ws.Range("J5").Value = 0.08265588920593574
$ws.Range("K5").Value = 22.4
$ws.Range("L5").Value = 0.04943721032884572
$ws.Range("M5").Value = 39.682
$ws.Range("N5").Value = 0.06498853586636096
$ws.Range("O5").Value = 1.771517857142857
$ws.Range("P5").Value = 39.6
$ws.Range("Q5").Value = 0.06485424172944644
$ws.Range("R5").Value = 1.767857142857143
$ws.Range("S5").Value = 0.08200000000000074
$ws.Range("T5").Value = 0.002066428103422225
$ws.Range("U5").Value = 251.5
$ws.Range("V5").Value = 0.4118899443170652
$ws.Range("W5").Value = 0.03696369636963696
$ws.Range("X5").Value = 0.06759593019751652
$ws.Range("Y5").Value = -0.03063223382787955
$ws.Range("Z5").Value = 0.9244316408202571
$ws.Range("AA5").Value = 0.07640971928210055
$ws.Range("AB5").Value = 0.06060576143334025
$ws.Range("AC5").Value = 0.01580395784876031
$ws.Range("AD5").Value = 104.8
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 104.8
$ws.Range("AG5").Value = -146.7
$ws.Range("AH5").Value = 0.1464914733016494
$ws.Range("AI5").Value = 0.1525695152132771
$ws.Range("AJ5").Value = -0.3162319465402026
$ws.Range("AK5").Value = -0.3369315571887919
$ws.Range("AL5").Value = 3.49
$ws.Range("AM5").Value = 3.49
$ws.Range("AN5").Value = 2.454332552693208
$ws.Range("AO5").Value = 12.12034383954155
$ws.Range("AP5").Value = -3.43559718969555
$ws.Range("AQ5").Value = 12.12034383954155

# Row 6
$ws.Range("D6").Value = -0.123
$ws.Range("E6").Value = -0.46
$ws.Range("F6").Value = 0.101
$ws.Range("G6").Value = 0.08080721828480346
$ws.Range("H6").Value = 0.08080721828480346
$ws.Range("I6").Value = 0.04907067429482253
$ws.Range("J6").Value = 0.04148532197062287
$ws.Range("K6").Value = 143
$ws.Range("L6").Value = 0.004070247345800245
$ws.Range("M6").Value = 1200
$ws.Range("N6").Value = 0.02510591558135886
$ws.Range("O6").Value = 8.391608391608392
$ws.Range("P6").Value = 1200
$ws.Range("Q6").Value = 0.02510591558135886
$ws.Range("R6").Value = 8.391608391608392
$ws.Range("U6").Value = 8384
$ws.Range("V6").Value = 0.1754066635284272
$ws.Range("W6").Value = 0.005720366103430619
$ws.Range("X6").Value = 0.07594425998795647
$ws.Range("Y6").Value = -0.07022389388452585
$ws.Range("Z6").Value = 0.9199480494996097
$ws.Range("AA6").Value = 0.03816434102973781
$ws.Range("AB6").Value = 0.0611297035486161
$ws.Range("AC6").Value = -0.0229653625188783
$ws.Range("AD6").Value = 17829
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 17829
$ws.Range("AG6").Value = 9445
$ws.Range("AH6").Value = 0.2716737903133642
$ws.Range("AI6").Value = 0.4801001723395089
$ws.Range("AJ6").Value = 0.1649997816307813
$ws.Range("AK6").Value = 0.3284988870339455
$ws.Range("AL6").Value = 386
$ws.Range("AM6").Value = 386
$ws.Range("AN6").Value = 8.865738438587767
$ws.Range("AO6").Value = 4.466321243523316
$ws.Range("AP6").Value = 4.696668324216808
$ws.Range("AQ6").Value = 4.466321243523316

# Row 7
$ws.Range("B7").Value = "Hansard Global Plc (LSE:HSD)"
$ws.Range("D7").Value = -0.133
$ws.Range("E7").Value = -0.213
$ws.Range("G7").Value = 0.1245283018867925
$ws.Range("H7").Value = 0.1245283018867925
$ws.Range("I7").Value = 0.1207547169811321
$ws.Range("J7").Value = 0.1156091551578811
$ws.Range("K7").Value = 5.57
$ws.Range("L7").Value = 0.08757861635220127
$ws.Range("M7").Value = 7.43
$ws.Range("N7").Value = 0.08898203592814372
$ws.Range("O7").Value = 1.333931777378815
$ws.Range("P7").Value = 7.43
$ws.Range("Q7").Value = 0.08898203592814372
$ws.Range("R7").Value = 1.333931777378815
$ws.Range("U7").Value = 49
$ws.Range("V7").Value = 0.5868263473053892
$ws.Range("W7").Value = 0.1609826589595376
$ws.Range("X7").Value = 0.06232255270229763
$ws.Range("Y7").Value = 0.09866010625723994
$ws.Range("Z7").Value = -3.854545454545454
$ws.Range("AA7").Value = -0.4456207435176507
$ws.Range("AB7").Value = 0.06051694865171065
$ws.Range("AC7").Value = -0.5061376921693614
$ws.Range("AD7").Value = 3.71
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 3.71
$ws.Range("AG7").Value = -45.29
$ws.Range("AH7").Value = 0.04254099300538929
$ws.Range("AI7").Value = 0.1036023457134879
$ws.Range("AJ7").Value = -1.185291808427113
$ws.Range("AK7").Value = 3.433661865049281
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0
$ws.Range("AN7").Value = 0.4684343434343434
$ws.Range("AP7").Value = -5.718434343434343

# Row 8
$ws.Range("B8").Value = "Legal & General Group Plc (LSE:LGEN)"
$ws.Range("D8").Value = 0.0354
$ws.Range("E8").Value = 0.0409
$ws.Range("F8").Value = 0.05690000000000001
$ws.Range("G8").Value = 0.171672798782035
$ws.Range("H8").Value = 0.171672798782035
$ws.Range("I8").Value = 0.1123181510614903
$ws.Range("J8").Value = 0.09403086484639799
$ws.Range("K8").Value = 1547.7
$ws.Range("L8").Value = 0.08181616341030196
$ws.Range("M8").Value = 1337.2
$ws.Range("N8").Value = 0.06198086621180657
$ws.Range("O8").Value = 0.8639917296633713
$ws.Range("P8").Value = 1297.6
$ws.Range("Q8").Value = 0.06014535746069415
$ws.Range("R8").Value = 0.8384053757188085
$ws.Range("S8").Value = 39.59999999999991
$ws.Range("T8").Value = 0.02961411905474119
$ws.Range("U8").Value = 26868.4
$ws.Range("V8").Value = 1.245383417383566
$ws.Range("W8").Value = 0.1390178836083391
$ws.Range("X8").Value = 0.07438955011683683
$ws.Range("Y8").Value = 0.06462833349150228
$ws.Range("Z8").Value = -15.50557377049179
$ws.Range("AA8").Value = -1.458002511578968
$ws.Range("AB8").Value = 0.06108083801922159
$ws.Range("AC8").Value = -1.519083349598189
$ws.Range("AD8").Value = 7238.4
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 7238.4
$ws.Range("AG8").Value = -19630
$ws.Range("AH8").Value = 0.251221679253665
$ws.Range("AI8").Value = 0.3824175824175824
$ws.Range("AJ8").Value = -10.09565932935609
$ws.Range("AK8").Value = 2.472167648984938
$ws.Range("AL8").Value = 279.8
$ws.Range("AM8").Value = 279.8
$ws.Range("AN8").Value = 3.347236994219653
$ws.Range("AO8").Value = 7.593638313080771
$ws.Range("AP8").Value = -9.077456647398844
$ws.Range("AQ8").Value = 7.593638313080771

# --- Clear cells that should become empty ---
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("AO7").ClearContents()
$ws.Range("AQ7").ClearContents()
$ws.Range("F7").ClearContents()

Write-Output "Update complete"